# Commit: "Change names from *img to img*"
# Rename the seven "*img" sheets to "img*" (prefix moved from the end to
# the front of each name), and move the active/selected tab from
# "holiday" (old activeTab index 8) to the last sheet, "eimg" -> "imge"
# (new activeTab index 16).

$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# Make the renamed "imge" sheet (previously "eimg", last sheet) the active tab.
$wb.Worksheets.Item("imge").Activate()
